# Add an "APELLIDOS" / "TELEFONO" column (D) to the employee-loading
# template on sheet "Hoja1": rename the old CONDUCTOR header (column C)
# to APELLIDOS, introduce a new column D headed TELEFONO, extend the
# title merge and the grey/yellow banner formatting over the new column,
# and reset the selection back to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- bring column D's formatting in line with column C (rows 1-3) ---
# (use Copy + PasteSpecial formats, since Range.Style only carries the
# named cell style and not the direct number/font/fill/border formatting
# actually used by this sheet)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)   # xlPasteFormats

# --- update the header texts ---
$ws.Range("C3").Value = "APELLIDOS"
$ws.Range("D3").Value = "TELEFONO"

# --- size the new column ---
$ws.Columns.Item(4).ColumnWidth = 21.6

# --- extend the title banner merge from A1:C2 to A1:D2 ---
$ws.Range("A1:D2").Merge()

# --- restore the selection to A4 (single cell, not A4:XFD11) ---
$ws.Range("A4").Select()
